# Automatische test-sync: 2025-06-22 19:09:50
#
# Appends two new incoming-mail log rows to the "Logs" sheet and refreshes
# the "Dashboard" pivot-style summary (category order/counts) to reflect
# them.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1. Append the two new log rows (38 and 39) to the Logs sheet.
# ---------------------------------------------------------------------

$newRows = @(
    @{
        Row = 38
        A = "Sollicitatie salesfunctie"
        B = "mailmind.test@zohomail.eu"
        C = "Hierbij mijn sollicitatie voor de salesfunctie. CV in bijlage."
        D = "Sollicitatie / Vacature"
        E = "Beste,`nBedankt voor je sollicitatie voor de salesfunctie. We waarderen je interesse in ons bedrijf. Ik zal je CV zorgvuldig bekijken en contact met je opnemen als er verdere stappen nodig zijn.`nMet vriendelijke groet,`n[Naam] E-mailassistent"
        F = "2025-06-22 19:09:11"
        G = "Ja"
    },
    @{
        Row = 39
        A = "Vragen over nieuwsbrief"
        B = "mailmind.test@zohomail.eu"
        C = "Wanneer wordt de volgende nieuwsbrief verstuurd?"
        D = "Afmelding / Nieuwsbrief"
        E = "Beste klant,`nDank voor je interesse in onze nieuwsbrief. De volgende nieuwsbrief staat gepland om aanstaande vrijdag te worden verstuurd. Houd dus je inbox in de gaten!`nMet vriendelijke groet,`n[Bedrijfsnaam]"
        F = "2025-06-22 19:09:14"
        G = "Ja"
    }
)

foreach ($r in $newRows) {
    $logs.Cells.Item($r.Row, 1).Value = $r.A
    $logs.Cells.Item($r.Row, 2).Value = $r.B
    $logs.Cells.Item($r.Row, 3).Value = $r.C
    $logs.Cells.Item($r.Row, 4).Value = $r.D
    $logs.Cells.Item($r.Row, 5).Value = $r.E
    $logs.Cells.Item($r.Row, 6).Value = $r.F
    $logs.Cells.Item($r.Row, 7).Value = $r.G
}

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting (Categorie / Beantwoord columns)
#    so the two new rows are covered too: D2:D37 -> D2:D39,
#    G2:G37 -> G2:G39.
# ---------------------------------------------------------------------

$catFormats = $logs.Range("D2:D37").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D39"))
}

$answeredFormats = $logs.Range("G2:G37").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G39"))
}

# ---------------------------------------------------------------------
# 3. Refresh the Dashboard summary table (rows 3-8) to reflect the two
#    new log entries: "Sollicitatie / Vacature" and
#    "Afmelding / Nieuwsbrief" both move from 3 to 4 occurrences, so the
#    category list is re-sorted descending by count.
# ---------------------------------------------------------------------

$dashboardRows = @(
    @{ Row = 3; Category = "Sollicitatie / Vacature"; Count = 4 },
    @{ Row = 4; Category = "Productinformatie"; Count = 4 },
    @{ Row = 5; Category = "Retour / Terugbetaling"; Count = 4 },
    @{ Row = 6; Category = "Afmelding / Nieuwsbrief"; Count = 4 },
    @{ Row = 7; Category = "Samenwerking / Partnerverzoek"; Count = 4 },
    @{ Row = 8; Category = "Offerte / Prijsaanvraag"; Count = 3 }
)

foreach ($r in $dashboardRows) {
    $dashboard.Cells.Item($r.Row, 1).Value = $r.Category
    $dashboard.Cells.Item($r.Row, 2).Value = $r.Count
}
